# Generate Report for Handback
#
# 1. The "Ready for handoff" status (shared across Overview!B3/C3 and the
#    zh-cn / de-de sheets' C3 cells) is now "Handback transform failed".
# 2. Each locale sheet (zh-cn, de-de) gets a new "Error Detail" (column L)
#    value on row 3 describing the handback/handoff file name mismatch.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 1. Update every cell that previously held "Ready for handoff" to the new
#    "Handback transform failed" status text so the shared string is fully
#    replaced (no stale references remain).
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# 2. Populate the "Error Detail" column (L) on row 3 for each locale sheet.
$wsZhCn.Range("L3").Value = "Handback file name: d0ktmxze.z1n is different with handoff file name: 74d5b395-6fa5-44f1-a494-c046316fcdaa.dd748f8f35e7966fdd3e3c7de6c2ee6451827450.zh-cn."
$wsDeDe.Range("L3").Value = "Handback file name: d0ktmxze.z1n is different with handoff file name: 74d5b395-6fa5-44f1-a494-c046316fcdaa.dd748f8f35e7966fdd3e3c7de6c2ee6451827450.de-de."
